$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 12: order status changes from Pending to Paid
$ws.Range("G12").Value = "Paid"

# Add new row 13
$ws.Range("A13").Value = "89bdc2f6-0e22-47a8-b4f2-b7b5696fc495"
$ws.Range("B13").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C13").Value = "P2001"
$ws.Range("D13").Value = "Approach"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 73.5
$ws.Range("G13").Value = "Delivered"
$ws.Range("H13").Value = "2025-08-07 23:02:17"

# Add new row 14
$ws.Range("A14").Value = "0947da20-6ab3-444d-97b4-2aa9c1662a75"
$ws.Range("B14").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C14").Value = "P2005"
$ws.Range("D14").Value = "Husband"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 182.9
$ws.Range("G14").Value = "Pending"
$ws.Range("H14").Value = "2025-08-07 23:14:23"
